$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the time_taken (F) column timestamps on the "data" sheet ---
# Rows 2..84 each get a refreshed timestamp (panel re-queried later the
# same day), in the same order as the sheet's existing rows.
$newTimestamps = @(
    "2021-10-05 14:20:27.031792",
    "2021-10-05 14:20:27.031800",
    "2021-10-05 14:20:27.031803",
    "2021-10-05 14:20:27.031806",
    "2021-10-05 14:20:27.031809",
    "2021-10-05 14:20:27.031812",
    "2021-10-05 14:20:27.031814",
    "2021-10-05 14:20:27.031817",
    "2021-10-05 14:20:27.031820",
    "2021-10-05 14:20:27.031823",
    "2021-10-05 14:20:27.031825",
    "2021-10-05 14:20:27.031828",
    "2021-10-05 14:20:27.031830",
    "2021-10-05 14:20:27.031833",
    "2021-10-05 14:20:27.031835",
    "2021-10-05 14:20:27.031838",
    "2021-10-05 14:20:27.031841",
    "2021-10-05 14:20:27.031844",
    "2021-10-05 14:20:27.031846",
    "2021-10-05 14:20:27.031849",
    "2021-10-05 14:20:27.031851",
    "2021-10-05 14:20:27.031854",
    "2021-10-05 14:20:27.031856",
    "2021-10-05 14:20:27.031859",
    "2021-10-05 14:20:27.031862",
    "2021-10-05 14:20:27.031865",
    "2021-10-05 14:20:27.031867",
    "2021-10-05 14:20:27.031870",
    "2021-10-05 14:20:27.031872",
    "2021-10-05 14:20:27.031875",
    "2021-10-05 14:20:27.031878",
    "2021-10-05 14:20:27.031880",
    "2021-10-05 14:20:27.031883",
    "2021-10-05 14:20:27.031886",
    "2021-10-05 14:20:27.031888",
    "2021-10-05 14:20:27.031891",
    "2021-10-05 14:20:27.031893",
    "2021-10-05 14:20:27.031896",
    "2021-10-05 14:20:27.031898",
    "2021-10-05 14:20:27.031901",
    "2021-10-05 14:20:27.031904",
    "2021-10-05 14:20:27.031907",
    "2021-10-05 14:20:27.031909",
    "2021-10-05 14:20:27.031912",
    "2021-10-05 14:20:27.031914",
    "2021-10-05 14:20:27.031917",
    "2021-10-05 14:20:27.031919",
    "2021-10-05 14:20:27.031922",
    "2021-10-05 14:20:27.031924",
    "2021-10-05 14:20:27.031927",
    "2021-10-05 14:20:27.031930",
    "2021-10-05 14:20:27.031932",
    "2021-10-05 14:20:27.031935",
    "2021-10-05 14:20:27.031938",
    "2021-10-05 14:20:27.031940",
    "2021-10-05 14:20:27.031943",
    "2021-10-05 14:20:27.031945",
    "2021-10-05 14:20:27.031948",
    "2021-10-05 14:20:27.031950",
    "2021-10-05 14:20:27.031953",
    "2021-10-05 14:20:27.031956",
    "2021-10-05 14:20:27.031958",
    "2021-10-05 14:20:27.031961",
    "2021-10-05 14:20:27.031964",
    "2021-10-05 14:20:27.031967",
    "2021-10-05 14:20:27.031970",
    "2021-10-05 14:20:27.031973",
    "2021-10-05 14:20:27.031975",
    "2021-10-05 14:20:27.031978",
    "2021-10-05 14:20:27.031981",
    "2021-10-05 14:20:27.031983",
    "2021-10-05 14:20:27.031986",
    "2021-10-05 14:20:27.031988",
    "2021-10-05 14:20:27.031991",
    "2021-10-05 14:20:27.031993",
    "2021-10-05 14:20:27.031996",
    "2021-10-05 14:20:27.032001",
    "2021-10-05 14:20:27.032004",
    "2021-10-05 14:20:27.032006",
    "2021-10-05 14:20:27.032009",
    "2021-10-05 14:20:27.032012",
    "2021-10-05 14:20:27.032014",
    "2021-10-05 14:20:27.032017"
)
for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = 2 + $i
    $dataSheet.Range("F$row").Value = $newTimestamps[$i]
}

# --- Add the new "metadata" tab, placed right after "data" ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row (B1:G1) - reuse the "data" sheet's bold/centered/bordered
# header style so the new tab matches the workbook's existing look.
$headers = @("data_name","data_id","data_version","data_version_created","panel_query_time","panel_get_request")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}
$dataSheet.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Data row 2: one summary record describing the "data" sheet's source panel.
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Gastrointestinal epithelial barrier disorders"
$ws.Range("C2").Value = 33
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.60"
$ws.Range("E2").Value = "2021-04-07T09:52:11.007055Z"
$ws.Range("F2").Value = "2021-10-05 14:20:27.028223"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/33/?format=json"

# D2 must stay text ("1.60") but carry no explicit cell style, matching its
# unstyled neighbours in row 2 - reset formatting from a never-touched,
# default-style cell without disturbing the stored text value.
$ws.Range("C3").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$dataSheet.Activate()
